# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (Home) updated with Week 17 totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 179
$wsOff.Range("C2").Value = 130
$wsOff.Range("D2").Value = 51
$wsOff.Range("E2").Value = 25

# --- DEF sheet: row 2 (Home) updated with Week 17 totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 333
$wsDef.Range("C2").Value = 231
$wsDef.Range("D2").Value = 59
$wsDef.Range("E2").Value = 27
$wsDef.Range("F2").Value = 4
$wsDef.Range("G2").Value = 3
